# [EXTRA SCRAPE] full data scraped for extra batting and bowling fields
#
# 1) Add a new "ODI Batting Extra" sheet (after "ODI Bowling") with the
#    scraped BATTING_POSITION / NUM_4 / NUM_6 / PERCENT_RUNS_OF_TOTAL /
#    MAN_OF_MATCH fields for each match.
# 2) Clear the (previously blank) INNING_NUMBER cells B3:B7 on "ODI Batting".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New sheet "ODI Batting Extra", placed after the last existing sheet
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Reuse the same bold/centered/bordered header style as the other sheets
# by copying the formatting (not the values) from "ODI Batting" row 1.
$battingWs = $wb.Worksheets.Item("ODI Batting")
$headerSrc = $battingWs.Range("A1:F1")
$headerSrc.Copy()
$headerDst = $ws.Range("A1:F1")
$headerDst.PasteSpecial(-4122)

$ws.Cells.Item(1, 1).Value = "MATCH_CODE"
$ws.Cells.Item(1, 2).Value = "BATTING_POSITION"
$ws.Cells.Item(1, 3).Value = "NUM_4"
$ws.Cells.Item(1, 4).Value = "NUM_6"
$ws.Cells.Item(1, 5).Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Cells.Item(1, 6).Value = "MAN_OF_MATCH"

# MATCH_CODE (column A) -- stored as text, one row per match
$matchCodes = @("4406", "4410", "4452", "4453", "4455", "4563", "4566", "4568")
$row = 2
foreach ($code in $matchCodes) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $code
    $row = $row + 1
}

# BATTING_POSITION (column B) -- numeric, only known for some rows
$ws.Cells.Item(3, 2).Value = 10
$ws.Cells.Item(4, 2).Value = 9
$ws.Cells.Item(5, 2).Value = 9
$ws.Cells.Item(6, 2).Value = 9
$ws.Cells.Item(7, 2).Value = 9
$ws.Cells.Item(8, 2).Value = 10

# NUM_4 (C), NUM_6 (D), PERCENT_RUNS_OF_TOTAL (E) -- only known for row 8,
# stored as text in the source data.
$c8 = $ws.Cells.Item(8, 3)
$c8.NumberFormat = "@"
$c8.Value = "0"

$d8 = $ws.Cells.Item(8, 4)
$d8.NumberFormat = "@"
$d8.Value = "0"

$e8 = $ws.Cells.Item(8, 5)
$e8.NumberFormat = "@"
$e8.Value = "1.14%"

# MAN_OF_MATCH (column F) -- text, "NO" for every data row
$row = 2
while ($row -le 9) {
    $ws.Cells.Item($row, 6).Value = "NO"
    $row = $row + 1
}

# ---------------------------------------------------------------------
# 2) Clear INNING_NUMBER (column B) for rows 3-7 on "ODI Batting"
# ---------------------------------------------------------------------
$battingWs.Range("B3:B7").ClearContents()

Write-Output "ODI Batting Extra sheet added and ODI Batting cleared."
